# "Generate Report for Handback" — mark the eae22f06 / 37d256ca entries as
# handed back (in sync with en-US), and populate the "Latest Target File" /
# "Latest Handback File" columns (E/F) with the handback hyperlinks, plus
# stamp the "Latest Handback DateTime" column (G) with the handback time.

$wb = $excel.ActiveWorkbook

# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    wherever it currently appears (Overview rollup + both language sheets).
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Text -eq "Ready for handoff") {
            $cell.Value = "Handed back: in sync with en-US"
        }
    }
}

# 2) Per-language-sheet handback details.
#    Row 2 -> 37d256ca-...md entry, Row 3 -> eae22f06-...md entry.
$sheetInfo = @{
    "zh-cn" = @{
        HandoffCommit  = "88bc9718bbd5e9815022c8e7af8b30b936fc7c9a"
        HandbackTime2  = "2016-02-22 05:22:46"
        HandbackTime3  = "2016-02-22 05:22:46"
    }
    "de-de" = @{
        HandoffCommit  = "7eca301f0928acf3588ca448604d0d09ddebe9cc"
        HandbackTime2  = "2016-02-22 05:23:08"
        HandbackTime3  = "2016-02-22 05:23:08"
    }
}

$rows = @(
    @{ Row = 2; MdFile = "37d256ca-dbbe-4435-b3c4-dee87b987bff.md"; XlfFile = "37d256ca-dbbe-4435-b3c4-dee87b987bff.43034988461c74100cfe5a4324f5d66ee11a0d47" },
    @{ Row = 3; MdFile = "eae22f06-c0d3-4ab5-8d24-39996e515c4f.md"; XlfFile = "eae22f06-c0d3-4ab5-8d24-39996e515c4f.6b95439b3453300fe3d7fb5100542f3f216c31c2" }
)

foreach ($langName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($langName)
    $info = $sheetInfo[$langName]

    foreach ($rowInfo in $rows) {
        $r = $rowInfo.Row
        $mdFile = $rowInfo.MdFile
        $xlfFileBase = $rowInfo.XlfFile
        $xlfFile = "$xlfFileBase.$langName.xlf"

        $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/7da99ec4bfc30c9d42b4a7f9831a7c7a44c45c76/e2e/$mdFile"
        $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($info.HandoffCommit)/ol-handoff/OpenLocalizationTestOrg/oltest.$langName/xinjiang/ht/$xlfFile"

        # E = Latest Target File (mirrors column A's source-file hyperlink)
        $ws.Hyperlinks.Add($ws.Range("E$r"), $mdUrl, "", "", $mdFile)
        # F = Latest Handback File (mirrors column C's handoff-file hyperlink)
        $ws.Hyperlinks.Add($ws.Range("F$r"), $xlfUrl, "", "", $xlfFile)

        # G = Latest Handback DateTime
        if ($r -eq 2) {
            $ws.Range("G$r").Value = $info.HandbackTime2
        } else {
            $ws.Range("G$r").Value = $info.HandbackTime3
        }
    }
}
